# "Modification on student enable" - rename the "Sub Class" header to
# "Campus" and move it to become the first column of the student sheet.
#
# Net column layout change (row 1 and all data rows):
#   old: A..N = Surname..Present Class, O = Sub Class, P..R = Home/Phone/Email
#   new: A = Campus, B..O = Surname..Present Class, P..R = Home/Phone/Email (unchanged)
#
# Achieved as: insert a blank column before A (shifts A..R right to B..S,
# so the old "Sub Class" column now sits at P), then delete that now-empty
# role at P (removing the duplicate "Sub Class"/O data that shifted there)
# so P..R collapse back onto the original Home/Phone/Email columns. Finally
# write "Campus" into the freshly inserted A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A; everything (A:R) shifts right to (B:S).
$ws.Columns("A:A").Insert(-4161)

# The old "Sub Class" header/column (originally O) is now at P. Delete it
# so the trailing Home Address/Phone Number/Email Address columns slide
# back down onto P:R (their original letters).
$ws.Columns("P:P").Delete(-4159)

# Rename "Sub Class" -> "Campus" as the new first column's header.
$ws.Range("A1").Value = "Campus"

# The leftover blank formatted cells that used to decorate column C
# (rows 3-21) are removed entirely rather than shifting to column D.
$ws.Range("D3:D21").Clear()
